$d = $word.ActiveDocument

# The document's header/footer images were inserted with their
# wp:docPr/@name (and pic:cNvPr/@name) pointing at the *other*
# image's default filename (image1.png <-> image2.png for the
# Pearson logo footers, image1.jpg <-> image2.jpg for the BTEC
# logo headers). This renames each inline picture's Name back to
# match its own media part, swapping "1" and "2" in the filename
# while keeping the extension untouched.

function Swap-ImageName($name) {
    if ($name -match "^(.*image)1(\.[A-Za-z]+)$") {
        return $matches[1] + "2" + $matches[2]
    } elseif ($name -match "^(.*image)2(\.[A-Za-z]+)$") {
        return $matches[1] + "1" + $matches[2]
    } else {
        return $name
    }
}

foreach ($sec in $d.Sections) {
    foreach ($hdr in $sec.Headers) {
        for ($i = 1; $i -le $hdr.Range.InlineShapes.Count; $i++) {
            $shp = $hdr.Range.InlineShapes.Item($i)
            $shp.Name = Swap-ImageName $shp.Name
        }
    }
    foreach ($ftr in $sec.Footers) {
        for ($i = 1; $i -le $ftr.Range.InlineShapes.Count; $i++) {
            $shp = $ftr.Range.InlineShapes.Item($i)
            $shp.Name = Swap-ImageName $shp.Name
        }
    }
}
